$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = 'sv'
$ws.Range("J14").Value = 'Statement-opinion'
$ws.Range("I29").Value = 'sv'
$ws.Range("J29").Value = 'Statement-opinion'
$ws.Range("I34").Value = 'aa'
$ws.Range("J34").Value = 'Agree/Accept'
$ws.Range("I44").Value = 'sd'
$ws.Range("J44").Value = 'Statement-non-opinion'
$ws.Range("I47").Value = 'sd'
$ws.Range("J47").Value = 'Statement-non-opinion'
$ws.Range("I52").Value = 'sd'
$ws.Range("J52").Value = 'Statement-non-opinion'
$ws.Range("I55").Value = 'sv'
$ws.Range("J55").Value = 'Statement-opinion'
$ws.Range("I83").Value = 'sv'
$ws.Range("J83").Value = 'Statement-opinion'
$ws.Range("I85").Value = 'sv'
$ws.Range("J85").Value = 'Statement-opinion'
$ws.Range("I101").Value = 'aa'
$ws.Range("J101").Value = 'Agree/Accept'
$ws.Range("I106").Value = '%'
$ws.Range("J106").Value = 'Uninterpretable'
$ws.Range("I110").Value = 'sv'
$ws.Range("J110").Value = 'Statement-opinion'
$ws.Range("I112").Value = 'sv'
$ws.Range("J112").Value = 'Statement-opinion'
$ws.Range("I114").Value = 'sv'
$ws.Range("J114").Value = 'Statement-opinion'
$ws.Range("I116").Value = 'b'
$ws.Range("J116").Value = 'Acknowledge (Backchannel)'
$ws.Range("I123").Value = 'sd'
$ws.Range("J123").Value = 'Statement-non-opinion'
$ws.Range("I125").Value = 'aa'
$ws.Range("J125").Value = 'Agree/Accept'
$ws.Range("I134").Value = 'aa'
$ws.Range("J134").Value = 'Agree/Accept'
$ws.Range("I135").Value = 'b'
$ws.Range("J135").Value = 'Acknowledge (Backchannel)'
$ws.Range("I136").Value = 'aa'
$ws.Range("J136").Value = 'Agree/Accept'
$ws.Range("I140").Value = '%'
$ws.Range("J140").Value = 'Uninterpretable'
$ws.Range("I141").Value = 'ba'
$ws.Range("J141").Value = 'Appreciation'
$ws.Range("I143").Value = 'sd'
$ws.Range("J143").Value = 'Statement-non-opinion'
$ws.Range("I147").Value = 'aa'
$ws.Range("J147").Value = 'Agree/Accept'
$ws.Range("I172").Value = 'sv'
$ws.Range("J172").Value = 'Statement-opinion'
$ws.Range("I175").Value = 'aa'
$ws.Range("J175").Value = 'Agree/Accept'
$ws.Range("I184").Value = 'qy'
$ws.Range("J184").Value = 'Yes-No-Question'
$ws.Range("I190").Value = 'aa'
$ws.Range("J190").Value = 'Agree/Accept'
$ws.Range("I204").Value = 'sd'
$ws.Range("J204").Value = 'Statement-non-opinion'
$ws.Range("I220").Value = 'sv'
$ws.Range("J220").Value = 'Statement-opinion'
$ws.Range("I237").Value = 'sv'
$ws.Range("J237").Value = 'Statement-opinion'
$ws.Range("I245").Value = 'aa'
$ws.Range("J245").Value = 'Agree/Accept'
$ws.Range("I249").Value = 'sd'
$ws.Range("J249").Value = 'Statement-non-opinion'
$ws.Range("I251").Value = 'sv'
$ws.Range("J251").Value = 'Statement-opinion'
$ws.Range("I253").Value = 'sv'
$ws.Range("J253").Value = 'Statement-opinion'
$ws.Range("I259").Value = 'aa'
$ws.Range("J259").Value = 'Agree/Accept'
$ws.Range("I260").Value = 'sd'
$ws.Range("J260").Value = 'Statement-non-opinion'
$ws.Range("I264").Value = 'ba'
$ws.Range("J264").Value = 'Appreciation'
$ws.Range("I265").Value = 'aa'
$ws.Range("J265").Value = 'Agree/Accept'
$ws.Range("I269").Value = 'aa'
$ws.Range("J269").Value = 'Agree/Accept'
$ws.Range("I281").Value = 'b'
$ws.Range("J281").Value = 'Acknowledge (Backchannel)'
$ws.Range("I286").Value = 'sd'
$ws.Range("J286").Value = 'Statement-non-opinion'
$ws.Range("I291").Value = 'sd'
$ws.Range("J291").Value = 'Statement-non-opinion'
$ws.Range("I305").Value = 'ba'
$ws.Range("J305").Value = 'Appreciation'
$ws.Range("I326").Value = 'sv'
$ws.Range("J326").Value = 'Statement-opinion'
$ws.Range("I352").Value = 'aa'
$ws.Range("J352").Value = 'Agree/Accept'
$ws.Range("I358").Value = 'sd'
$ws.Range("J358").Value = 'Statement-non-opinion'
$ws.Range("I362").Value = 'sd'
$ws.Range("J362").Value = 'Statement-non-opinion'
$ws.Range("I364").Value = 'b'
$ws.Range("J364").Value = 'Acknowledge (Backchannel)'
$ws.Range("I368").Value = 'sd'
$ws.Range("J368").Value = 'Statement-non-opinion'
$ws.Range("I370").Value = 'aa'
$ws.Range("J370").Value = 'Agree/Accept'
$ws.Range("I375").Value = 'sv'
$ws.Range("J375").Value = 'Statement-opinion'
$ws.Range("I386").Value = '%'
$ws.Range("J386").Value = 'Uninterpretable'
$ws.Range("I387").Value = 'aa'
$ws.Range("J387").Value = 'Agree/Accept'
$ws.Range("I388").Value = 'sd'
$ws.Range("J388").Value = 'Statement-non-opinion'
$ws.Range("I391").Value = 'ba'
$ws.Range("J391").Value = 'Appreciation'
$ws.Range("I392").Value = 'sd'
$ws.Range("J392").Value = 'Statement-non-opinion'
$ws.Range("I408").Value = 'b'
$ws.Range("J408").Value = 'Acknowledge (Backchannel)'
$ws.Range("I414").Value = 'sd'
$ws.Range("J414").Value = 'Statement-non-opinion'
$ws.Range("I430").Value = 'aa'
$ws.Range("J430").Value = 'Agree/Accept'
$ws.Range("I441").Value = 'sd'
$ws.Range("J441").Value = 'Statement-non-opinion'
$ws.Range("I456").Value = 'sv'
$ws.Range("J456").Value = 'Statement-opinion'
$ws.Range("I458").Value = 'sv'
$ws.Range("J458").Value = 'Statement-opinion'
$ws.Range("I464").Value = 'ba'
$ws.Range("J464").Value = 'Appreciation'
$ws.Range("I470").Value = 'sd'
$ws.Range("J470").Value = 'Statement-non-opinion'
$ws.Range("I472").Value = 'sv'
$ws.Range("J472").Value = 'Statement-opinion'
$ws.Range("I477").Value = '%'
$ws.Range("J477").Value = 'Uninterpretable'
$ws.Range("I480").Value = 'sd'
$ws.Range("J480").Value = 'Statement-non-opinion'
$ws.Range("I494").Value = 'sv'
$ws.Range("J494").Value = 'Statement-opinion'
$ws.Range("I497").Value = 'sd'
$ws.Range("J497").Value = 'Statement-non-opinion'
$ws.Range("I499").Value = 'sv'
$ws.Range("J499").Value = 'Statement-opinion'
$ws.Range("I508").Value = 'sd'
$ws.Range("J508").Value = 'Statement-non-opinion'
$ws.Range("I517").Value = 'b'
$ws.Range("J517").Value = 'Acknowledge (Backchannel)'
$ws.Range("I521").Value = 'aa'
$ws.Range("J521").Value = 'Agree/Accept'
$ws.Range("I530").Value = 'sv'
$ws.Range("J530").Value = 'Statement-opinion'
